# "50%" sheet - finish the "alles wat ik wil upgraden" line and add the
# still-missing answer for the buck/boost converter / line sensor question.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F4 was empty - fill in the answer about the main/line sensor improvements.
$ws.Range("F4").Value = "ik weet hoe ik het hooft en de lijnsensoren wil verbeteren bij het hoofdbord zal ik allen nog de buck/boost covertor onderdelen opzoeken"

# C10 had a typo ("wiel" -> "wil"); fix the text.
$ws.Range("C10").Value = "alles wat ik wil upgraden"

# Update the view: scroll down a bit and select C10:F10 (the cell that was
# just edited), matching where the author's cursor ended up when saving.
$ws.Range("C10:F10").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 1
